$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.465.50'
$ws.Range('E2').Value = '  -0.88%  '
$ws.Range('D3').Value = '1.625.56'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.75'
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.500'
$ws.Range('E6').Value = '  +1.42%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('E9').Value = '  -2.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.68'
$ws.Range('E10').Value = '  -1.87%  '
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('D12').Value = '1.851.86'
$ws.Range('E12').Value = '  -0.77%  '
$ws.Range('D13').Value = '1.638.71'
$ws.Range('E13').Value = '  -0.07%  '
$ws.Range('E14').Value = '  +1.28%  '
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.83'
$ws.Range('E16').Value = '  +2.87%  '
$ws.Range('D17').Value = '26.505.49'
$ws.Range('E17').Value = '  -0.78%  '
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '213.57'
$ws.Range('E19').Value = '  +2.53%  '
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('E22').Value = '  +1.46%  '
$ws.Range('E23').Value = '  -1.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.05'
$ws.Range('E24').Value = '  +6.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.48'
$ws.Range('E25').Value = '  +1.31%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E27').Value = '  -0.93%  '
$ws.Range('E28').Value = '  +0.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.47'
$ws.Range('E29').Value = '  +0.50%  '
$ws.Range('E30').Value = '  -1.66%  '
$ws.Range('E31').Value = '  -0.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.31'
$ws.Range('E32').Value = '  +2.57%  '
$ws.Range('E33').Value = '  -0.76%  '
$ws.Range('E34').Value = '  -0.40%  '
$ws.Range('D35').Value = '1.224.69'
$ws.Range('E35').Value = '  +4.90%  '
$ws.Range('E36').Value = '  -1.71%  '
$ws.Range('E37').Value = '  +3.24%  '
$ws.Range('E38').Value = '  +0.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.792'
$ws.Range('E39').Value = '  -2.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.504'
$ws.Range('E40').Value = '  +0.18%  '
$ws.Range('E41').Value = '  -1.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.792'
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.32'
$ws.Range('E43').Value = '  -0.95%  '
$ws.Range('D44').Value = '1.761.64'
$ws.Range('E44').Value = '  -0.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.67'
$ws.Range('E45').Value = '  +0.16%  '
$ws.Range('E46').Value = '  +0.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.68'
$ws.Range('E47').Value = '  -0.49%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0510'
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.407'
$ws.Range('E49').Value = '  -0.88%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.49'
$ws.Range('E50').Value = '  -0.47%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.27%  '
